# remove simulationName #415 (#441)
#
# The "simulationName" input row (and its accompanying "display name of
# simulation for report" description) is removed from the SimulationSets
# sheet. Deleting the entire row shifts every row below it up by one and
# keeps all per-row formatting/styles/validations attached to their
# (renamed) content, matching how Excel itself performs a row delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimulationSets")

# Row 4 holds: A4 = "simulationName", B4 = "display name of simulation for report"
$ws.Rows.Item(4).Delete()

# Reflect the author's final UI state: cursor left on the SimulationSets
# sheet, sitting on the row that used to be "outputs" (now row 4), with
# cell A5 (the next row, "observedDataFile") selected.
$ws.Activate()
$ws.Range("A5").Select()
